# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" sheet (with fund-holding detail data) right after the
# "总计" (summary) sheet, shifting the existing "2022-Q2" / "2022-Q1" /
# "2021-Q4" sheets one position to the right, and updates the "总计" summary
# sheet with a new row for 2022-Q3 (and re-numbers the existing rows).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert the 2022-Q3 row at the top
#    of the data and push the other rows down by one.
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryRows = @(
    @("2022-Q3", 34, 7.87),
    @("2022-Q2", 11, 2.98),
    @("2022-Q1", 25, 5.93),
    @("2021-Q4", 4, 0.1)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]

    $idxCell = $summary.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# -----------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计" and rename it.
# -----------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Match the page-margin conventions used by the sibling sheets in this workbook.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q3.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$dataText = @"
0|720001|财通价值动量混合|38.35|79.13|7.61|2.9184|1
1|001480|财通成长优选混合|20.31|91.20|7.38|1.4989|2
2|005106|银华农业产业股票A|12.00|93.82|5.21|0.6252|8
3|014915|财通匠心优选一年持有期混合A|5.65|81.89|7.51|0.4243|2
4|008983|财通科技创新混合A|2.95|94.28|7.95|0.2345|1
5|164403|前海开源沪港深农业混合（LOF）A|3.43|89.98|6.14|0.2106|5
6|005421|中欧嘉泽灵活配置混合|6.99|83.96|3.01|0.2104|10
7|013993|中欧光熠一年持有期混合型证券投资基金A|6.06|84.26|3.35|0.2030|9
8|011708|中欧嘉益一年混合A|4.44|93.35|3.96|0.1758|8
9|010418|财通景气行业混合A|2.72|94.88|5.79|0.1575|9
10|009062|财通智慧成长混合A|2.17|84.78|6.19|0.1343|6
11|015210|前海开源沪港深农业混合（LOF）C|1.93|89.98|6.14|0.1185|5
12|009063|财通智慧成长混合C|1.50|84.78|6.19|0.0928|6
13|008984|财通科技创新混合C|1.16|94.28|7.95|0.0922|1
14|501015|财通多策略升级混合（LOF）A|2.06|94.80|4.41|0.0908|9
15|013994|中欧光熠一年持有期混合型证券投资基金C|2.45|84.26|3.35|0.0821|9
16|005270|太平改革红利精选灵活配置混合|1.78|82.55|4.30|0.0765|9
17|011709|中欧嘉益一年混合C|1.87|93.35|3.96|0.0741|8
18|010897|太平价值增长股票C|0.91|83.64|6.10|0.0555|8
19|001940|农银汇理现代农业加灵活配置混合|1.06|74.84|4.90|0.0519|3
20|005959|财通新视野灵活配置混合C|1.12|94.59|4.56|0.0511|9
21|519678|银河消费驱动混合A|0.96|91.11|5.04|0.0484|7
22|014916|财通匠心优选一年持有期混合C|0.61|81.89|7.51|0.0458|2
23|002844|金鹰多元策略灵活配置混合|0.48|89.55|8.73|0.0419|2
24|014064|银华农业产业股票C|0.75|93.82|5.21|0.0391|8
25|005851|财通新视野灵活配置混合A|0.62|94.59|4.56|0.0283|9
26|015271|财通多策略升级混合（LOF）C|0.62|94.80|4.41|0.0273|9
27|350009|天治研究驱动混合A|0.29|93.73|8.69|0.0252|2
28|002043|天治研究驱动混合C|0.24|93.73|8.69|0.0209|2
29|010896|太平价值增长股票A|0.15|83.64|6.10|0.0092|8
30|015041|西部利得新富灵活配置混合C|0.12|68.35|4.65|0.0056|3
31|673120|西部利得新富灵活配置混合A|0.03|68.35|4.65|0.0014|3
32|015668|银河消费驱动混合C|0.01|91.11|5.04|0.0005|7
33|016234|财通景气行业混合C|0.00|94.88|5.79|0|9
"@

$lines = $dataText -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $f = $line -split "\|"

    $r = [int]$f[0] + 2

    $idxCell = $q3.Cells.Item($r, 1)
    $idxCell.Value = [int]$f[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $codeCell = $q3.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $f[1]

    $q3.Cells.Item($r, 3).Value = $f[2]

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $f[3]

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $f[4]

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $f[5]

    if ($f[6] -eq "0") {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        $q3.Cells.Item($r, 7).NumberFormat = "@"
        $q3.Cells.Item($r, 7).Value = $f[6]
    }

    $q3.Cells.Item($r, 8).Value = [int]$f[7]
}

# Restore the originally active sheet/selection.
$summary.Activate()
$null = $summary.Range("A1").Select()

